$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 142.94118
$ws.Range("I39").Value = 142.94118
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 428.82354
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -132.82354
$ws.Range("N39").ClearContents()
$ws.Range("H51").Value = 4998.6665
$ws.Range("J51").Value = 4998.6665
$ws.Range("L51").Value = 4998.6665
$ws.Range("N51").Value = -5966.6665
$ws.Range("H86").Value = 19825
$ws.Range("I86").Value = 32500
$ws.Range("J86").Value = 7150
$ws.Range("K86").Value = 32500
$ws.Range("L86").Value = 7150
$ws.Range("M86").Value = -31377
$ws.Range("N86").Value = -9396
$ws.Range("H89").Value = 19825
$ws.Range("I89").Value = 32500
$ws.Range("J89").Value = 7150
$ws.Range("K89").Value = 162500
$ws.Range("L89").Value = 35750
$ws.Range("M89").Value = -156884
$ws.Range("N89").Value = -46982
$ws.Range("H113").Value = 4491.6665
$ws.Range("I113").Value = 4490
$ws.Range("K113").Value = 4490
$ws.Range("M113").Value = -1236

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8889.793
$ws.Range("I32").Value = 7421.5713
$ws.Range("K32").Value = 7421.5713
$ws.Range("M32").Value = -7134.5713
$ws.Range("H132").Value = 2242.24
$ws.Range("I132").Value = 2002.7142
$ws.Range("K132").Value = 6008.142599999999
$ws.Range("M132").Value = -3478.142599999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 21986.75
$ws.Range("I75").Value = 6556.2856
$ws.Range("J75").Value = 130000
$ws.Range("K75").Value = 6556.2856
$ws.Range("L75").Value = 130000
$ws.Range("M75").Value = -5620.2856
$ws.Range("N75").Value = -131872
$ws.Range("H78").Value = 21986.75
$ws.Range("I78").Value = 6556.2856
$ws.Range("J78").Value = 130000
$ws.Range("K78").Value = 19668.8568
$ws.Range("L78").Value = 390000
$ws.Range("M78").Value = -14988.8568
$ws.Range("N78").Value = -399360
$ws.Range("H94").Value = 1135.625
$ws.Range("I94").Value = 1132.3914
$ws.Range("K94").Value = 1132.3914
$ws.Range("M94").Value = -681.3914
$ws.Range("H134").Value = 5895.6924
$ws.Range("I134").Value = 6404
$ws.Range("K134").Value = 19212
$ws.Range("M134").Value = -16677

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 49666.668
$ws.Range("J64").Value = 49666.668
$ws.Range("L64").Value = 49666.668
$ws.Range("N64").Value = -50162.668
$ws.Range("H67").Value = 49666.668
$ws.Range("J67").Value = 49666.668
$ws.Range("L67").Value = 49666.668
$ws.Range("N67").Value = -51382.668
$ws.Range("H93").Value = 88203.5
$ws.Range("I93").Value = 88203.5
$ws.Range("K93").Value = 88203.5
$ws.Range("M93").Value = -86331.5
$ws.Range("H103").Value = 6883
$ws.Range("I103").Value = 6883
$ws.Range("K103").Value = 6883
$ws.Range("M103").Value = -5711
$ws.Range("H105").Value = 909.2727
$ws.Range("I105").Value = 909.2727
$ws.Range("K105").Value = 909.2727
$ws.Range("M105").Value = 837.7273

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 194.41667
$ws.Range("J12").Value = 207.54546
$ws.Range("L12").Value = 622.6363799999999
$ws.Range("N12").Value = -968.6363799999999
$ws.Range("H127").Value = 26666.666
$ws.Range("J127").Value = 26666.666
$ws.Range("L127").Value = 79999.99800000001
$ws.Range("N127").Value = -89919.99800000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 63000
$ws.Range("J15").Value = 63000
$ws.Range("L15").Value = 63000
$ws.Range("N15").Value = -63576
$ws.Range("H47").Value = 19000
$ws.Range("I47").Value = 7000
$ws.Range("K47").Value = 7000
$ws.Range("M47").Value = -6432
$ws.Range("H55").Value = 27295.166
$ws.Range("I55").Value = 6997.5
$ws.Range("J55").Value = 37444
$ws.Range("K55").Value = 6997.5
$ws.Range("L55").Value = 37444
$ws.Range("M55").Value = -6670.5
$ws.Range("N55").Value = -38098
$ws.Range("H81").Value = 63000
$ws.Range("J81").Value = 63000
$ws.Range("L81").Value = 63000
$ws.Range("N81").Value = -64996
$ws.Range("H84").Value = 63000
$ws.Range("J84").Value = 63000
$ws.Range("L84").Value = 189000
$ws.Range("N84").Value = -198984
$ws.Range("H97").Value = 672.3913
$ws.Range("I97").Value = 660.2381
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 660.2381
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -164.2381
$ws.Range("N97").Value = -1792

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H46").Value = 4499
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4499
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4499
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4875
$ws.Range("H55").Value = 758.9
$ws.Range("I55").Value = 1380
$ws.Range("J55").Value = 344.83334
$ws.Range("K55").Value = 1380
$ws.Range("L55").Value = 344.83334
$ws.Range("M55").Value = -1207
$ws.Range("N55").Value = -690.83334
$ws.Range("H68").Value = 7000
$ws.Range("I68").Value = 7000
$ws.Range("K68").Value = 7000
$ws.Range("M68").Value = -6251
$ws.Range("H71").Value = 7000
$ws.Range("I71").Value = 7000
$ws.Range("K71").Value = 35000
$ws.Range("M71").Value = -31256
$ws.Range("H100").Value = 2181.5
$ws.Range("I100").Value = 2181.5
$ws.Range("K100").Value = 2181.5
$ws.Range("M100").Value = -1640.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 41000
$ws.Range("J52").Value = 41000
$ws.Range("L52").Value = 41000
$ws.Range("N52").Value = -41452
$ws.Range("H126").Value = 3708.3076
$ws.Range("I126").Value = 4521.6
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 13564.8
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -11094.8
$ws.Range("N126").Value = -14540
$ws.Range("H132").Value = 815.6667
$ws.Range("I132").Value = 739.8
$ws.Range("K132").Value = 2219.4
$ws.Range("M132").Value = 310.6000000000004
